# Handback status report generation: add a new handed-back file entry
# ("ea5545c1-0e13-4200-9bdd-d9e193e46182.md") alongside the existing
# ("88f50d06-952f-460a-936e-51cbecf97898.md", renamed from
# "6bbff5fe-2d23-4f51-baa8-e64361157706.md") row, on every sheet / table.

$wb = $excel.ActiveWorkbook

$HYPER_UNDERLINE = 2        # xlUnderlineStyleSingle
$HYPER_COLOR = 15570276     # OLE BGR for RGB(0x64,0x95,0xED) -> matches existing HyperLink font
$DATE_FMT = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $HYPER_UNDERLINE
    $rng.Font.Color = $HYPER_COLOR
}

function Style-AsDate($rng) {
    $rng.NumberFormat = $DATE_FMT
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Row 2: rename the existing handed-back file from 6bbff5fe... to 88f50d06...
$wsOv.Range("A2").Value = "88f50d06-952f-460a-936e-51cbecf97898.md"
$wsOv.Range("B2").Value = "e2e\88f50d06-952f-460a-936e-51cbecf97898.md"
$wsOv.Range("C2").Value = ".md"
$wsOv.Range("E2").Value = "Handed back: in sync with en-US"
$wsOv.Range("F2").Value = "Handed back: in sync with en-US"
$wsOv.Range("G2").Value = "2016-08-15 10:58:37"

# Row 3: new handed-back file entry
$wsOv.Range("A3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$wsOv.Range("B3").Value = "e2e\ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-15 10:58:37"

Style-AsDate($wsOv.Range("G2"))
Style-AsDate($wsOv.Range("G3"))

# Hyperlinks: drop the old link(s) entirely and rebuild rId2 (renamed row) + rId3 (new row)
$wsOv.Range("B2:B3").Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/88f50d06-952f-460a-936e-51cbecf97898.md", "", "", "e2e\88f50d06-952f-460a-936e-51cbecf97898.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/ea5545c1-0e13-4200-9bdd-d9e193e46182.md", "", "", "e2e\ea5545c1-0e13-4200-9bdd-d9e193e46182.md")
Style-AsHyperlink($wsOv.Range("B2"))
Style-AsHyperlink($wsOv.Range("B3"))

$loOv = $wsOv.ListObjects.Item("Overview")
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 2: rename handed-back file, refresh handoff/handback timestamps
$ws2.Range("A2").Value = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "e2e"
$ws2.Range("E2").Value = "ht"
$ws2.Range("F2").Value = "False"
$ws2.Range("G2").Value = "88f50d06-952f-460a-936e-51cbecf97898.ac5fd329c46a4ef799f14df0faf2a074131fc266.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-15 10:58:30"
$ws2.Range("I2").Value = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws2.Range("J2").Value = "88f50d06-952f-460a-936e-51cbecf97898.ac5fd329c46a4ef799f14df0faf2a074131fc266.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-15 10:58:57"
$ws2.Range("L2").Value = ""
$ws2.Range("M2").Value = "True"
$ws2.Range("N2").Value = ""
$ws2.Range("O2").Value = "False"
$ws2.Range("P2").Value = ""

# Row 3: new handed-back file entry (content duplicate of row 2)
$ws2.Range("A3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.c3e6c39f554c15d98b5d943ccf97ee479f27c8b0.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-15 10:58:30"
$ws2.Range("I3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$ws2.Range("J3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.c3e6c39f554c15d98b5d943ccf97ee479f27c8b0.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-15 10:58:57"
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

Style-AsDate($ws2.Range("H2"))
Style-AsDate($ws2.Range("K2"))
Style-AsDate($ws2.Range("H3"))
Style-AsDate($ws2.Range("K3"))

# Hyperlinks: repoint existing A2/I2, add new A3/I3
$ws2.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/88f50d06-952f-460a-936e-51cbecf97898.md"
$ws2.Hyperlinks.Item(1).TextToDisplay = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws2.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/927ac11912f0fd37b4613bd80ba3378be373c8bc/e2e/88f50d06-952f-460a-936e-51cbecf97898.md"
$ws2.Hyperlinks.Item(2).TextToDisplay = "88f50d06-952f-460a-936e-51cbecf97898.md"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/ea5545c1-0e13-4200-9bdd-d9e193e46182.md", "", "", "ea5545c1-0e13-4200-9bdd-d9e193e46182.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/927ac11912f0fd37b4613bd80ba3378be373c8bc/e2e/ea5545c1-0e13-4200-9bdd-d9e193e46182.md", "", "", "ea5545c1-0e13-4200-9bdd-d9e193e46182.md")

Style-AsHyperlink($ws2.Range("A2"))
Style-AsHyperlink($ws2.Range("I2"))
Style-AsHyperlink($ws2.Range("A3"))
Style-AsHyperlink($ws2.Range("I3"))

$lo2 = $ws2.ListObjects.Item("zh-cn")
$lo2.Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 2: rename handed-back file, refresh handoff/handback timestamps
$ws3.Range("A2").Value = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "e2e"
$ws3.Range("E2").Value = "ht"
$ws3.Range("F2").Value = "False"
$ws3.Range("G2").Value = "88f50d06-952f-460a-936e-51cbecf97898.ac5fd329c46a4ef799f14df0faf2a074131fc266.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-15 10:58:37"
$ws3.Range("I2").Value = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws3.Range("J2").Value = "88f50d06-952f-460a-936e-51cbecf97898.ac5fd329c46a4ef799f14df0faf2a074131fc266.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-15 10:59:09"
$ws3.Range("L2").Value = ""
$ws3.Range("M2").Value = "True"
$ws3.Range("N2").Value = ""
$ws3.Range("O2").Value = "False"
$ws3.Range("P2").Value = ""

# Row 3: new handed-back file entry (content duplicate of row 2)
$ws3.Range("A3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.c3e6c39f554c15d98b5d943ccf97ee479f27c8b0.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-15 10:58:37"
$ws3.Range("I3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.md"
$ws3.Range("J3").Value = "ea5545c1-0e13-4200-9bdd-d9e193e46182.c3e6c39f554c15d98b5d943ccf97ee479f27c8b0.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-15 10:59:09"
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

Style-AsDate($ws3.Range("H2"))
Style-AsDate($ws3.Range("K2"))
Style-AsDate($ws3.Range("H3"))
Style-AsDate($ws3.Range("K3"))

# Hyperlinks: repoint existing A2/I2, add new A3/I3
$ws3.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/88f50d06-952f-460a-936e-51cbecf97898.md"
$ws3.Hyperlinks.Item(1).TextToDisplay = "88f50d06-952f-460a-936e-51cbecf97898.md"
$ws3.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ba685b1bd0aba2765bd53d42f1d20509c10ccb2f/e2e/88f50d06-952f-460a-936e-51cbecf97898.md"
$ws3.Hyperlinks.Item(2).TextToDisplay = "88f50d06-952f-460a-936e-51cbecf97898.md"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/ea5545c1-0e13-4200-9bdd-d9e193e46182.md", "", "", "ea5545c1-0e13-4200-9bdd-d9e193e46182.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ba685b1bd0aba2765bd53d42f1d20509c10ccb2f/e2e/ea5545c1-0e13-4200-9bdd-d9e193e46182.md", "", "", "ea5545c1-0e13-4200-9bdd-d9e193e46182.md")

Style-AsHyperlink($ws3.Range("A2"))
Style-AsHyperlink($ws3.Range("I2"))
Style-AsHyperlink($ws3.Range("A3"))
Style-AsHyperlink($ws3.Range("I3"))

$lo3 = $ws3.ListObjects.Item("de-de")
$lo3.Resize($ws3.Range("A1:P3"))
